$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet 1: "Neg_Change"
# Remove rows 3-6 (KOTAKBANK, EICHERMOT, LT, ITC), leaving only the
# header row and row 2, then overwrite row 2 with the BEL data.
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("A3:A6").EntireRow.Delete()

$ws1.Cells.Item(2,1).Value = "BEL"
$ws1.Cells.Item(2,2).Value = 420
$ws1.Cells.Item(2,3).Value = 422.4
$ws1.Cells.Item(2,4).Value = 414.5
$ws1.Cells.Item(2,5).Value = 415.75
$ws1.Cells.Item(2,6).Value = 12831686
$ws1.Cells.Item(2,7).Value = 27222212
$ws1.Cells.Item(2,8).Value = -0.5286317658535611
$ws1.Cells.Item(2,9).Value = "BEL"

# ---------------------------------------------------------------
# Sheet 2: "Pos_Change"
# Update rows 2-3 (SHRIRAMFIN -> BAJFINANCE, AXISBANK -> SUNPHARMA)
# and append rows 4-5 (INFY, EICHERMOT).
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Cells.Item(2,1).Value = "BAJFINANCE"
$ws2.Cells.Item(2,2).Value = 1040
$ws2.Cells.Item(2,3).Value = 1060.4
$ws2.Cells.Item(2,4).Value = 1035
$ws2.Cells.Item(2,5).Value = 1054.7
$ws2.Cells.Item(2,6).Value = 6202540
$ws2.Cells.Item(2,7).Value = 4135617
$ws2.Cells.Item(2,8).Value = 0.4997858844278859
$ws2.Cells.Item(2,9).Value = "BAJFINANCE"

$ws2.Cells.Item(3,1).Value = "SUNPHARMA"
$ws2.Cells.Item(3,2).Value = 1707
$ws2.Cells.Item(3,3).Value = 1711.2
$ws2.Cells.Item(3,4).Value = 1682.2
$ws2.Cells.Item(3,5).Value = 1702.7
$ws2.Cells.Item(3,6).Value = 1597993
$ws2.Cells.Item(3,7).Value = 1099955
$ws2.Cells.Item(3,8).Value = 0.4527803410139506
$ws2.Cells.Item(3,9).Value = "SUNPHARMA"

$ws2.Cells.Item(4,1).Value = "INFY"
$ws2.Cells.Item(4,2).Value = 1479.7
$ws2.Cells.Item(4,3).Value = 1481.9
$ws2.Cells.Item(4,4).Value = 1462.9
$ws2.Cells.Item(4,5).Value = 1468
$ws2.Cells.Item(4,6).Value = 8691330
$ws2.Cells.Item(4,7).Value = 5470600
$ws2.Cells.Item(4,8).Value = 0.5887343253025262
$ws2.Cells.Item(4,9).Value = "INFY"

$ws2.Cells.Item(5,1).Value = "EICHERMOT"
$ws2.Cells.Item(5,2).Value = 7050
$ws2.Cells.Item(5,3).Value = 7060
$ws2.Cells.Item(5,4).Value = 6906.5
$ws2.Cells.Item(5,5).Value = 6917
$ws2.Cells.Item(5,6).Value = 494329
$ws2.Cells.Item(5,7).Value = 334962
$ws2.Cells.Item(5,8).Value = 0.4757763567210609
$ws2.Cells.Item(5,9).Value = "EICHERMOT"

Write-Host "edit complete"
